$wb = $excel.ActiveWorkbook

# --- Rename sheets to new timestamped task-order names ---
$wb.Worksheets.Item("GNG_TO-165029112392062").Name = "GNG_TO-16504777813653793"
$wb.Worksheets.Item("NB_TO-16502911260204015").Name = "NB_TO-16504777823934112"
$wb.Worksheets.Item("RS_TO-16502911260214045").Name = "RS_TO-16504777823943772"
$wb.Worksheets.Item("TOL_TO-1650291126078502").Name = "TOL_TO-1650477782457409"
$wb.Worksheets.Item("vSAT_TO-16502911261397064").Name = "vSAT_TO-16504777825193763"

# --- GNG sheet (sheet1): update stimulus file names in column B ---
$wsGNG = $wb.Worksheets.Item("GNG_TO-16504777813653793")
$wsGNG.Range("B2").Value = "go_stims-1650477781324376.csv"
$wsGNG.Range("B3").Value = "GNG_stims-16504777813493748.csv"
$wsGNG.Range("B4").Value = "go_stims-16504777813503752.csv"
$wsGNG.Range("B5").Value = "GNG_stims-16504777813644102.csv"

# --- NB sheet (sheet2): update stimulus file names in column B ---
$wsNB = $wb.Worksheets.Item("NB_TO-16504777823934112")
$wsNB.Range("B2").Value = "ZB-match_3-16504777814363751.csv"
$wsNB.Range("B3").Value = "OB-1650477781675409.csv"
$wsNB.Range("B4").Value = "ZB-match_6-16504777814764092.csv"
$wsNB.Range("B5").Value = "TB-165047778237041.csv"
$wsNB.Range("B6").Value = "TB-1650477782344375.csv"
$wsNB.Range("B7").Value = "OB-16504777820424118.csv"
$wsNB.Range("B8").Value = "TB-16504777822754128.csv"
$wsNB.Range("B9").Value = "OB-16504777821644094.csv"
$wsNB.Range("B10").Value = "ZB-match_6-165047778139538.csv"

# --- TOL sheet (sheet4): update stimulus file names in column B ---
$wsTOL = $wb.Worksheets.Item("TOL_TO-1650477782457409")
$wsTOL.Range("B2").Value = "MM_stims-1650477782409375.csv"
$wsTOL.Range("B3").Value = "ZM_stims-16504777823963804.csv"
$wsTOL.Range("B4").Value = "MM_stims-16504777824413745.csv"
$wsTOL.Range("B5").Value = "ZM_stims-1650477782410375.csv"
$wsTOL.Range("B6").Value = "MM_stims-1650477782457409.csv"
$wsTOL.Range("B7").Value = "ZM_stims-16504777824423785.csv"

# --- vSAT sheet (sheet5): update stimulus file names in column B ---
$wsVSAT = $wb.Worksheets.Item("vSAT_TO-16504777825193763")
$wsVSAT.Range("B2").Value = "vSAT_stims-16504777824873755.csv"
$wsVSAT.Range("B3").Value = "SAT_stims-16504777824603775.csv"
$wsVSAT.Range("B4").Value = "vSAT_stims-16504777825033758.csv"
$wsVSAT.Range("B5").Value = "SAT_stims-16504777824723778.csv"
